$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "C5"  = 164
    "C6"  = 199
    "C8"  = 337
    "I11" = 1940
    "C12" = 695
    "C13" = 779
    "C14" = 892
    "C15" = 887
    "C16" = 1064
    "C17" = 1006
    "C18" = 1114
    "C19" = 1219
    "C20" = 1302
    "C21" = 1491
    "C23" = 1957
    "I24" = 10384
    "C25" = 2498
    "C26" = 2761
    "I26" = 13874
    "C27" = 3184
    "I27" = 16382
    "C29" = 4188
    "C30" = 4745
    "I30" = 21817
    "C31" = 5114
    "I31" = 24857
    "C32" = 5623
    "I32" = 28296
    "I33" = 29140
    "C34" = 6179
    "I34" = 32394
    "C35" = 6448
    "I35" = 35508
    "C36" = 7373
    "I36" = 42631
    "C37" = 11875
    "I37" = 67359
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
